$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Legend")

# 1. Rename the data sheet from "Sheet1" to "Data".
#    Excel automatically updates the _xlnm._FilterDatabase defined name
#    (and any other references) to point at the new sheet name.
$ws1.Name = "Data"

# 2. Add a header row to the top of the Legend sheet with "Column1" / "Column2".
$ws2.Range("A1").EntireRow.Insert()
$ws2.Range("A1").Value = "Column1"
$ws2.Range("B1").Value = "Column2"

# 3. Convert the Legend range into a proper Excel Table named "Table1".
$tbl = $ws2.ListObjects.Add(1, $ws2.Range("A1:B7"), $null, 1)
$tbl.Name = "Table1"

# Approximate the manually-set width of the new "Column2" column.
$ws2.Columns.Item(2).ColumnWidth = 10.46

# 4. Make the Legend sheet the active tab, with A1:B7 selected (matches the
#    tabSelected flag moving from the Data sheet to the Legend sheet).
$ws2.Activate()
$ws2.Range("A1:B7").Select() | Out-Null
